$d = $word.ActiveDocument

$pairs = @(
    @("879÷6=", "274÷2="),
    @("726÷9=", "687÷2="),
    @("770÷5=", "777÷7="),
    @("417÷8=", "786÷9="),
    @("414÷5=", "569÷7="),
    @("640÷5=", "876÷3="),
    @("580÷2=", "435÷9="),
    @("105÷3=", "970÷7="),
    @("520÷4=", "742÷4="),
    @("834÷9=", "860÷6="),
    @("956÷4=", "359÷3="),
    @("558÷3=", "979÷6="),
    @("542÷3=", "950÷5="),
    @("836÷3=", "620÷4="),
    @("763÷4=", "585÷4="),
    @("522÷9=", "553÷3="),
    @("763÷9=", "349÷4="),
    @("166÷9=", "568÷2="),
    @("696÷3=", "346÷8="),
    @("614÷2=", "298÷8="),
    @("787÷8=", "808÷8="),
    @("434÷6=", "833÷9="),
    @("801÷3=", "943÷4="),
    @("740÷6=", "961÷4="),
    @("412÷6=", "773÷5=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
